$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 - Department of Health and Human Services: current_count 43 -> 44, change 0 -> 1
$ws.Range("B10").Value = 44
$ws.Range("D10").Value = 1

# Row 26 - Intelligence Community: current_count 25 -> 26, change 0 -> 1
$ws.Range("B26").Value = 26
$ws.Range("D26").Value = 1
